$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("settings")
$wsSurvey   = $wb.Worksheets.Item("survey")
$wsModel    = $wb.Worksheets.Item("model")

# --- survey sheet: insert a new row before the old last data row and
#     fill it with the new "AMOSTRA" / "integer" entry -------------------
$wsSurvey.Activate()
$null = $wsSurvey.Rows.Item(10).Insert()
$wsSurvey.Range("D10").Value = "integer"
$wsSurvey.Range("F10").Value = "AMOSTRA"
$null = $wsSurvey.Range("F10").Select()

# --- model sheet: append a new row for "AMOSTRA" -------------------------
$wsModel.Activate()
$wsModel.Range("A9").Value = "AMOSTRA"
$wsModel.Range("B9").Value = "integer"
$wsModel.Range("C9").Value = $false
$null = $wsModel.Range("A10").Select()

# --- settings sheet becomes the active / selected tab --------------------
$wsSettings.Activate()
